$d = $word.ActiveDocument

for ($i = 1; $i -le 6; $i++) {
    $old = "<id>p088r_a$i</id>"
    $new = "<id>p088r_$i</id>"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Host "done"
